$wb = $excel.ActiveWorkbook
$wsEin = $wb.Worksheets.Item("Einnahmen")
$wsAus = $wb.Worksheets.Item("Ausgaben")

# ---------------------------------------------------------------------------
# 1) Ausgaben (expenses): invoice for "Film: Garfield" -> new row 105
#    (formats copied from the last existing data row, 104, so the new row
#    reuses the same style indices instead of creating new ones)
# ---------------------------------------------------------------------------
$wsAus.Range("A104:K104").Copy()
$wsAus.Range("A105:K105").PasteSpecial(-4122)

$wsAus.Range("A105").Value = "Verleiher"
$wsAus.Range("B105").Value = 45543
$wsAus.Range("C105").Value = "Film: Garfield"
$wsAus.Range("D105").Value = 45614
$wsAus.Range("E105").Value = 86.75
$wsAus.Range("F105").Value = "Sony Pictures Releasing Switzerland GmbH"
$wsAus.Range("G105").Value = "Dufourstrasse 59, 8008 Zürich"
$wsAus.Range("H105").Value = "275862000000000000030626311"
$wsAus.Range("I105").Value = "3062631"
$wsAus.Range("J105").Value = "4404"
$wsAus.Range("K105").Value = "Filmmiete Kino"

# ---------------------------------------------------------------------------
# 2) Einnahmen (income): Förderer (sponsor) amounts from Theater am Bahnhof
#    -> new rows 19-21
# ---------------------------------------------------------------------------
$wsEin.Range("A17:E17").Copy()
$wsEin.Range("A19:E19").PasteSpecial(-4122)
$wsEin.Range("A20:E20").PasteSpecial(-4122)
$wsEin.Range("A21:E21").PasteSpecial(-4122)

$wsEin.Range("A19").Value = "Sonstiges"
$wsEin.Range("B19").Value = "Förderer Kino"
$wsEin.Range("C19").Value = 45657
$wsEin.Range("D19").Value = 12800
$wsEin.Range("E19").Value = "Theater am Bahnhof"

$wsEin.Range("A20").Value = "Sonstiges"
$wsEin.Range("B20").Value = "Förderer Firmen Kino"
$wsEin.Range("C20").Value = 45657
$wsEin.Range("D20").Value = 1750
$wsEin.Range("E20").Value = "Theater am Bahnhof"

$wsEin.Range("A21").Value = "Sonstiges"
$wsEin.Range("B21").Value = "Förderer TaB nur Kino-Anteil"
$wsEin.Range("C21").Value = 45657
$wsEin.Range("D21").Value = 3306
$wsEin.Range("E21").Value = "Theater am Bahnhof"

# grow the Einnahmen table ("Table3") to cover the 3 new rows
$tblEin = $wsEin.ListObjects.Item("Table3")
$tblEin.Resize($wsEin.Range("A1:G21"))

# ---------------------------------------------------------------------------
# 3) Ausgaben (expenses): invoice for "Film: Tabubrecherin" -> new row 106
# ---------------------------------------------------------------------------
$wsAus.Range("A104:K104").Copy()
$wsAus.Range("A106:K106").PasteSpecial(-4122)

$wsAus.Range("A106").Value = "Verleiher"
$wsAus.Range("B106").Value = 45634
$wsAus.Range("C106").Value = "Film: Tabubrecherin"
$wsAus.Range("D106").Value = 45665
$wsAus.Range("E106").Value = 172.25
$wsAus.Range("F106").Value = "langjahr film GmbH"
$wsAus.Range("G106").Value = "Bahnhofstrasse 7, Zug"
$wsAus.Range("J106").Value = "4404"
$wsAus.Range("K106").Value = "Filmmiete Kino"

# grow the Ausgaben table ("Table16") to cover the 2 new rows
$tblAus = $wsAus.ListObjects.Item("Table16")
$tblAus.Resize($wsAus.Range("A1:K106"))

# ---------------------------------------------------------------------------
# 4) Final view / selection state
# ---------------------------------------------------------------------------
$wsEin.Range("B22").Select()

$wsAus.Activate()
$wsAus.Range("A107").Select()
